# Edit the "model selection" sheet: split "XGBoost Regressor" into
# "XGBoost Regressor default" / "XGBoost Regressor tuned", add Random Forest
# results, add a new "trained/tested" flag column (E) and a "COR" metric
# column, then wrap the range in a real Excel Table (Table2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model selection")

# ---- Header row (row 5) ----
$ws.Range("D5").Value = "Model"
$ws.Range("E5").Value = "MAE"
$ws.Range("F5").Value = "RMSE"
$ws.Range("G5").Value = "COR"
$ws.Range("H5").Value = "R² (R-Squared)"
$ws.Range("I5").Value = "Training Time"
$ws.Range("J5").Value = "Inference Time"
$ws.Range("K5").Value = "Interpretability"
$ws.Range("L5").Value = "Robustness to Outliers"

# ---- Data rows (D6:L14) ----
# Row 6: XGBoost Regressor tuned
$ws.Range("D6").Value = "XGBoost Regressor tuned"
$ws.Range("E6").Value = 1
$ws.Range("K6").Value = "Medium"
$ws.Range("L6").Value = "High"

# Row 7: XGBoost Regressor default (new backup GPU version)
$ws.Range("D7").Value = "XGBoost Regressor default"
$ws.Range("E7").Value = 1
$ws.Range("K7").Value = "Medium"
$ws.Range("L7").Value = "High"

# Row 8: Linear Regression
$ws.Range("D8").Value = "Linear Regression"
$ws.Range("E8").Value = 0
$ws.Range("K8").Value = "High"
$ws.Range("L8").Value = "Low"

# Row 9: Ridge Regression
$ws.Range("D9").Value = "Ridge Regression"
$ws.Range("K9").Value = "High"
$ws.Range("L9").Value = "Medium"

# Row 10: Lasso Regression
$ws.Range("D10").Value = "Lasso Regression"
$ws.Range("K10").Value = "High"
$ws.Range("L10").Value = "Medium"

# Row 11: Random Forest Regressor (newly trained + tested)
$ws.Range("D11").Value = "Random Forest Regressor"
$ws.Range("E11").Value = 1
$ws.Range("K11").Value = "Medium"
$ws.Range("L11").Value = "High"

# Row 12: Gradient Boosting Regressor
$ws.Range("D12").Value = "Gradient Boosting Regressor"
$ws.Range("K12").Value = "Medium"
$ws.Range("L12").Value = "High"

# Row 13: Support Vector Regressor
$ws.Range("D13").Value = "Support Vector Regressor"
$ws.Range("E13").Value = 0
$ws.Range("K13").Value = "Low"
$ws.Range("L13").Value = "Medium"

# Row 14: Neural Network Regressor
$ws.Range("D14").Value = "Neural Network Regressor"
$ws.Range("E14").Value = 1
$ws.Range("K14").Value = "Low"
$ws.Range("L14").Value = "Medium"

# ---- Column widths ----
$ws.Columns.Item(5).ColumnWidth = 7.33203125
$ws.Columns.Item(6).ColumnWidth = 7.33203125
$ws.Columns.Item(7).ColumnWidth = 7.33203125
$ws.Columns.Item(8).ColumnWidth = 15.33203125
$ws.Columns.Item(9).ColumnWidth = 15.33203125
$ws.Columns.Item(10).ColumnWidth = 15.33203125
$ws.Columns.Item(11).ColumnWidth = 14.5
$ws.Columns.Item(12).ColumnWidth = 20.83203125

# ---- Turn the range into a real table (Table2) ----
$tableRange = $ws.Range("D5:L14")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table2"
$lo.TableStyle = "TableStyleMedium2"

# ---- View state ----
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("A12:XFD12").Select()
